# Update the new n channel mosfet
# The four old N-channel MOSFETs (Q1, Q2, Q3, Q4 - part BSH103,235 /
# SOT95P230X110-3N) are replaced by four new ones (Q5, Q6, Q7, Q8 - part
# 2N7002K-7 / SOT96P240X100-3N), with the designator->coordinate pairing
# reshuffled across the four placement rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ----- Sheet1 (placement / Layer-Rotation sheet), rows 22-25 -----
$ws1.Range("A22").Value = "Q5"
$ws1.Range("B22").Value = -72.39
$ws1.Range("C22").Value = -48.26
$ws1.Range("E22").Value = 270

$ws1.Range("A23").Value = "Q6"
$ws1.Range("B23").Value = -80.01
$ws1.Range("C23").Value = -48.26
$ws1.Range("E23").Value = 270

$ws1.Range("A24").Value = "Q7"
$ws1.Range("B24").Value = -74.93
$ws1.Range("C24").Value = 38.1
$ws1.Range("E24").Value = 90

$ws1.Range("A25").Value = "Q8"
$ws1.Range("B25").Value = -64.77
$ws1.Range("C25").Value = -48.26
$ws1.Range("E25").Value = 270

# ----- Sheet2 (BOM sheet with Value/Footprint columns), rows 21-24 -----
$ws2.Range("A21").Value = "Q5"
$ws2.Range("B21").Value = -72.39
$ws2.Range("C21").Value = -48.26
$ws2.Range("D21").Value = 270
$ws2.Range("E21").Value = "2N7002K-7"
$ws2.Range("F21").Value = "SOT96P240X100-3N"

$ws2.Range("A22").Value = "Q6"
$ws2.Range("B22").Value = -80.01
$ws2.Range("C22").Value = -48.26
$ws2.Range("D22").Value = 270
$ws2.Range("E22").Value = "2N7002K-7"
$ws2.Range("F22").Value = "SOT96P240X100-3N"

$ws2.Range("A23").Value = "Q7"
$ws2.Range("B23").Value = -74.93
$ws2.Range("C23").Value = 38.1
$ws2.Range("D23").Value = 90
$ws2.Range("E23").Value = "2N7002K-7"
$ws2.Range("F23").Value = "SOT96P240X100-3N"

$ws2.Range("A24").Value = "Q8"
$ws2.Range("B24").Value = -64.77
$ws2.Range("C24").Value = -48.26
$ws2.Range("D24").Value = 270
$ws2.Range("E24").Value = "2N7002K-7"
$ws2.Range("F24").Value = "SOT96P240X100-3N"

# ----- Active-cell selection moved from L6 to H8 on Sheet1 -----
$ws1.Activate()
$ws1.Range("H8").Select()
